$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1113.6666
$ws.Range("I19").Value = 687.1
$ws.Range("J19").Value = 1364.5883
$ws.Range("K19").Value = 687.1
$ws.Range("L19").Value = 1364.5883
$ws.Range("M19").Value = -512.1
$ws.Range("N19").Value = -1714.5883
$ws.Range("H64").Value = 41837.848
$ws.Range("J64").Value = 3549.2727
$ws.Range("L64").Value = 3549.2727
$ws.Range("N64").Value = -4045.2727
$ws.Range("H67").Value = 41837.848
$ws.Range("J67").Value = 3549.2727
$ws.Range("L67").Value = 3549.2727
$ws.Range("N67").Value = -5265.2727
$ws.Range("H70").Value = 1860
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 1950
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 5850
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -6390
$ws.Range("H73").Value = 1860
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 1950
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 5850
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -7722
$ws.Range("H96").Value = 588.8125
$ws.Range("I96").Value = 351.2
$ws.Range("K96").Value = 1053.6
$ws.Range("M96").Value = 319.4000000000001
$ws.Range("H103").Value = 719
$ws.Range("I103").Value = 448.75
$ws.Range("J103").Value = 770.4761999999999
$ws.Range("K103").Value = 1346.25
$ws.Range("L103").Value = 2311.4286
$ws.Range("M103").Value = -760.25
$ws.Range("N103").Value = -3483.4286
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 86233.336
$ws.Range("I102").Value = 144934.28
$ws.Range("J102").Value = 4052
$ws.Range("K102").Value = 144934.28
$ws.Range("L102").Value = 4052
$ws.Range("M102").Value = -143312.28
$ws.Range("N102").Value = -7296
$ws.Range("H122").Value = 1664.75
$ws.Range("I122").Value = 1720
$ws.Range("K122").Value = 5160
$ws.Range("M122").Value = -2710
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1525.75
$ws.Range("I99").Value = 1487.1428
$ws.Range("J99").Value = 1579.8
$ws.Range("K99").Value = 1487.1428
$ws.Range("L99").Value = 1579.8
$ws.Range("M99").Value = 10.85719999999992
$ws.Range("N99").Value = -4575.8
$ws.Range("H105").Value = 155957.77
$ws.Range("I105").Value = 126660
$ws.Range("J105").Value = 202834.2
$ws.Range("K105").Value = 126660
$ws.Range("L105").Value = 202834.2
$ws.Range("M105").Value = -124913
$ws.Range("N105").Value = -206328.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23327.486
$ws.Range("I31").Value = 1004.0732
$ws.Range("K31").Value = 1004.0732
$ws.Range("M31").Value = -709.0732
$ws.Range("H34").Value = 23327.486
$ws.Range("I34").Value = 1004.0732
$ws.Range("K34").Value = 1004.0732
$ws.Range("M34").Value = -802.0732
$ws.Range("H62").Value = 2500
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 2500
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 447.5
$ws.Range("I86").Value = 447.5
$ws.Range("K86").Value = 1342.5
$ws.Range("M86").Value = -156.5
$ws.Range("H89").Value = 447.5
$ws.Range("I89").Value = 447.5
$ws.Range("K89").Value = 4027.5
$ws.Range("M89").Value = 1900.5
$ws.Range("H98").Value = 140000.12
$ws.Range("I98").Value = 1003
$ws.Range("J98").Value = 159856.86
$ws.Range("K98").Value = 3009
$ws.Range("L98").Value = 479570.58
$ws.Range("M98").Value = -1511
$ws.Range("N98").Value = -482566.58
$ws.Range("H131").Value = 831.8200000000001
$ws.Range("J131").Value = 866.54346
$ws.Range("L131").Value = 2599.63038
$ws.Range("N131").Value = -12679.63038
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 166668700
$ws.Range("I80").Value = 250002000
$ws.Range("J80").Value = 2100
$ws.Range("K80").Value = 250002000
$ws.Range("L80").Value = 2100
$ws.Range("M80").Value = -250001002
$ws.Range("N80").Value = -4096
$ws.Range("H83").Value = 166668700
$ws.Range("I83").Value = 250002000
$ws.Range("J83").Value = 2100
$ws.Range("K83").Value = 1250010000
$ws.Range("L83").Value = 10500
$ws.Range("M83").Value = -1250005008
$ws.Range("N83").Value = -20484
$ws.Range("H102").Value = 195020.73
$ws.Range("I102").Value = 1431.7084
$ws.Range("J102").Value = 858754.5600000001
$ws.Range("K102").Value = 1431.7084
$ws.Range("L102").Value = 858754.5600000001
$ws.Range("M102").Value = 190.2916
$ws.Range("N102").Value = -861998.5600000001
$ws.Range("H107").Value = 555.1667
$ws.Range("I107").Value = 459.08334
$ws.Range("J107").Value = 747.3333
$ws.Range("K107").Value = 459.08334
$ws.Range("L107").Value = 747.3333
$ws.Range("M107").Value = 1460.91666
$ws.Range("N107").Value = -4587.3333
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 73606.21000000001
$ws.Range("I40").Value = 201635.8
$ws.Range("J40").Value = 2478.6667
$ws.Range("K40").Value = 201635.8
$ws.Range("L40").Value = 2478.6667
$ws.Range("M40").Value = -201499.8
$ws.Range("N40").Value = -2750.6667
$ws.Range("H68").Value = 3276.3845
$ws.Range("I68").Value = 1633.6666
$ws.Range("J68").Value = 4684.4287
$ws.Range("K68").Value = 1633.6666
$ws.Range("L68").Value = 4684.4287
$ws.Range("M68").Value = -884.6666
$ws.Range("N68").Value = -6182.4287
$ws.Range("H71").Value = 3276.3845
$ws.Range("I71").Value = 1633.6666
$ws.Range("J71").Value = 4684.4287
$ws.Range("K71").Value = 8168.333000000001
$ws.Range("L71").Value = 23422.1435
$ws.Range("M71").Value = -4424.333000000001
$ws.Range("N71").Value = -30910.1435
$ws.Range("H82").Value = 1042.1538
$ws.Range("I82").Value = 811.2
$ws.Range("J82").Value = 1186.5
$ws.Range("K82").Value = 811.2
$ws.Range("L82").Value = 1186.5
$ws.Range("M82").Value = -450.2
$ws.Range("N82").Value = -1908.5
$ws.Range("H85").Value = 1042.1538
$ws.Range("I85").Value = 811.2
$ws.Range("J85").Value = 1186.5
$ws.Range("K85").Value = 811.2
$ws.Range("L85").Value = 1186.5
$ws.Range("M85").Value = 436.8
$ws.Range("N85").Value = -3682.5
$ws.Range("H100").Value = 1303.4286
$ws.Range("I100").Value = 1050
$ws.Range("K100").Value = 1050
$ws.Range("M100").Value = -509
$ws.Range("H132").Value = 2801.186
$ws.Range("I132").Value = 2786.639
$ws.Range("J132").Value = 2876
$ws.Range("K132").Value = 8359.917000000001
$ws.Range("L132").Value = 8628
$ws.Range("M132").Value = -5829.917000000001
$ws.Range("N132").Value = -13688
$ws.Range("H136").Value = 2234.5
$ws.Range("I136").Value = 2223.7778
$ws.Range("J136").Value = 2266.6667
$ws.Range("K136").Value = 6671.3334
$ws.Range("L136").Value = 6800.000100000001
$ws.Range("M136").Value = -4121.3334
$ws.Range("N136").Value = -11900.0001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 142858750
$ws.Range("I96").Value = 200001860
$ws.Range("J96").Value = 1002
$ws.Range("K96").Value = 200001860
$ws.Range("L96").Value = 1002
$ws.Range("M96").Value = -200000487
$ws.Range("N96").Value = -3748
$ws.Range("H126").Value = 1247.8422
$ws.Range("I126").Value = 1157.7858
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 3473.3574
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -1003.3574
$ws.Range("N126").Value = -9440
$ws.Range("H132").Value = 5104.4116
$ws.Range("I132").Value = 5750.2
$ws.Range("K132").Value = 17250.6
$ws.Range("M132").Value = -14720.6
